$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("D1").Value = "Memory Usage (bytes)"

# Update Run Time (ms) values in column C
$ws.Range("C2").Value = 17.84706115722656
$ws.Range("C3").Value = 17.25196838378906
$ws.Range("C4").Value = 17.44198799133301
$ws.Range("C5").Value = 20.07317543029785
$ws.Range("C6").Value = 17.85826683044434
